# "planilha de teste atualizada" - update the test guide worksheet:
#  - widen column C to fit the new observation text
#  - center-align (horizontal + vertical) the whole data table
#  - mark the "IOS" requirement row as not-met and explain why, wrapping the
#    long justification text and growing the row height to fit it
#  - move the active selection like the author left it

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- widen column C so the new justification text has room -----------------
# (64.6328125 is the stored character-width we're after; the host quantizes
# ColumnWidth writes to its own pixel grid, so 63.75 is the input that lands
# closest to it.)
$ws.Columns.Item(3).ColumnWidth = 63.75

# --- center align (horizontal + vertical) the body of the table -------------
# Each call below targets a single "border group" so every cell in the group
# picks up exactly one new combined style (no stray half-updated styles left
# behind from writing Horizontal/Vertical alignment separately on a big,
# heterogeneously-styled range).

# Row 2 (top row of the body, one cell per border group)
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4108

# Rows 3-11 (middle rows share one border style per column) - format one
# cell, then fan that exact format out to the rest of the column with a
# format-only paste so the shared borders/fill survive untouched.
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").Copy()
$ws.Range("A3:A11").PasteSpecial(-4122)

$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").Copy()
$ws.Range("B3:B11").PasteSpecial(-4122)

$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108
$ws.Range("C3").Copy()
$ws.Range("C3:C11").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- row 12 ("IOS"): mark not attended + add the justification text --------
$ws.Range("B12").Value = "Não atendido"
$ws.Range("C12").Value = "Por termos uma liguagem hibrida tinhamos o interesse de fazer os testes nas duas plataformas. Com o tempo de teste curto, o tempo para a criação de ambiente em um emulador de IOS e a adaptação do códova nesse ambiente nos custaria um tempo de dedicação maior, optamos por deixar o sistema somente em android, visto que ali a aplicação estava estavel. "

$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").VerticalAlignment = -4108

$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("B12").VerticalAlignment = -4108

$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("C12").VerticalAlignment = -4108
$ws.Range("C12").WrapText = $true

$ws.Rows.Item(12).RowHeight = 73

# --- selection, as left by the author ---------------------------------------
$ws.Range("C13").Select()
